# Auto-generated edit script applying numeric updates from the commit diff.
# Updates currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ
# columns (H-N) on specific rows across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: H62=3957.6956, I62=3497.1, J62=4312, K62=3497.1, L62=4312, M62=-2873.1, N62=-5560
$ws.Range("H62").Value = 3957.6956
$ws.Range("I62").Value = 3497.1
$ws.Range("J62").Value = 4312
$ws.Range("K62").Value = 3497.1
$ws.Range("L62").Value = 4312
$ws.Range("M62").Value = -2873.1
$ws.Range("N62").Value = -5560
# Row 65: H65=3957.6956, I65=3497.1, J65=4312, K65=17485.5, L65=21560, M65=-14365.5, N65=-27800
$ws.Range("H65").Value = 3957.6956
$ws.Range("I65").Value = 3497.1
$ws.Range("J65").Value = 4312
$ws.Range("K65").Value = 17485.5
$ws.Range("L65").Value = 21560
$ws.Range("M65").Value = -14365.5
$ws.Range("N65").Value = -27800
# Row 80: H80=18743586, I80=433.33334, J80=24366532, K80=1300.00002, L80=73099596, M80=-302.0000199999999, N80=-73101592
$ws.Range("H80").Value = 18743586
$ws.Range("I80").Value = 433.33334
$ws.Range("J80").Value = 24366532
$ws.Range("K80").Value = 1300.00002
$ws.Range("L80").Value = 73099596
$ws.Range("M80").Value = -302.0000199999999
$ws.Range("N80").Value = -73101592
# Row 83: H83=18743586, I83=433.33334, J83=24366532, K83=3900.00006, L83=219298788, M83=1091.99994, N83=-219308772
$ws.Range("H83").Value = 18743586
$ws.Range("I83").Value = 433.33334
$ws.Range("J83").Value = 24366532
$ws.Range("K83").Value = 3900.00006
$ws.Range("L83").Value = 219298788
$ws.Range("M83").Value = 1091.99994
$ws.Range("N83").Value = -219308772
# Row 129: H129=2356.0344, J129=2604.923, L129=7814.768999999999, N129=-17814.769
$ws.Range("H129").Value = 2356.0344
$ws.Range("J129").Value = 2604.923
$ws.Range("L129").Value = 7814.768999999999
$ws.Range("N129").Value = -17814.769

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32=7488.602, I32=5099.5244, J32=25298.092, K32=5099.5244, L32=25298.092, M32=-4812.5244, N32=-25872.092
$ws.Range("H32").Value = 7488.602
$ws.Range("I32").Value = 5099.5244
$ws.Range("J32").Value = 25298.092
$ws.Range("K32").Value = 5099.5244
$ws.Range("L32").Value = 25298.092
$ws.Range("M32").Value = -4812.5244
$ws.Range("N32").Value = -25872.092
# Row 74: H74=17544986, I74=23256340, K74=23256340, M74=-23255466
$ws.Range("H74").Value = 17544986
$ws.Range("I74").Value = 23256340
$ws.Range("K74").Value = 23256340
$ws.Range("M74").Value = -23255466
# Row 77: H77=17544986, I77=23256340, K77=116281700, M77=-116277332
$ws.Range("H77").Value = 17544986
$ws.Range("I77").Value = 23256340
$ws.Range("K77").Value = 116281700
$ws.Range("M77").Value = -116277332
# Row 97: H97=760.7, I97=769.8889, K97=769.8889, M97=-273.8889
$ws.Range("H97").Value = 760.7
$ws.Range("I97").Value = 769.8889
$ws.Range("K97").Value = 769.8889
$ws.Range("M97").Value = -273.8889
# Row 109: H109=19088.5, J109=19088.5, L109=19088.5, N109=-21862.5
$ws.Range("H109").Value = 19088.5
$ws.Range("J109").Value = 19088.5
$ws.Range("L109").Value = 19088.5
$ws.Range("N109").Value = -21862.5
# Row 122: H122=1668.4546, I122=1546.2821, K122=4638.846299999999, M122=-2188.846299999999
$ws.Range("H122").Value = 1668.4546
$ws.Range("I122").Value = 1546.2821
$ws.Range("K122").Value = 4638.846299999999
$ws.Range("M122").Value = -2188.846299999999
# Row 124: H124=26497.5, J124=26497.5, L124=26497.5, N124=-36317.5
$ws.Range("H124").Value = 26497.5
$ws.Range("J124").Value = 26497.5
$ws.Range("L124").Value = 26497.5
$ws.Range("N124").Value = -36317.5
# Row 125: H125=30861.334, J125=30861.334, L125=30861.334, N125=-40701.334
$ws.Range("H125").Value = 30861.334
$ws.Range("J125").Value = 30861.334
$ws.Range("L125").Value = 30861.334
$ws.Range("N125").Value = -40701.334
# Row 132: H132=11962.902, I132=1871.35, K132=5614.049999999999, M132=-3084.049999999999
$ws.Range("H132").Value = 11962.902
$ws.Range("I132").Value = 1871.35
$ws.Range("K132").Value = 5614.049999999999
$ws.Range("M132").Value = -3084.049999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94: H94=1176.0625, I94=952.25, J94=1399.875, K94=952.25, L94=1399.875, M94=-501.25, N94=-2301.875
$ws.Range("H94").Value = 1176.0625
$ws.Range("I94").Value = 952.25
$ws.Range("J94").Value = 1399.875
$ws.Range("K94").Value = 952.25
$ws.Range("L94").Value = 1399.875
$ws.Range("M94").Value = -501.25
$ws.Range("N94").Value = -2301.875
# Row 107: H107=1406.0465, I107=1095.2778, J107=3004.2856, K107=1095.2778, L107=3004.2856, M107=824.7221999999999, N107=-6844.2856
$ws.Range("H107").Value = 1406.0465
$ws.Range("I107").Value = 1095.2778
$ws.Range("J107").Value = 3004.2856
$ws.Range("K107").Value = 1095.2778
$ws.Range("L107").Value = 3004.2856
$ws.Range("M107").Value = 824.7221999999999
$ws.Range("N107").Value = -6844.2856
# Row 134: H134=2926.5576, I134=3049.3125, J134=1453.5, K134=9147.9375, L134=4360.5, M134=-6612.9375, N134=-9430.5
$ws.Range("H134").Value = 2926.5576
$ws.Range("I134").Value = 3049.3125
$ws.Range("J134").Value = 1453.5
$ws.Range("K134").Value = 9147.9375
$ws.Range("L134").Value = 4360.5
$ws.Range("M134").Value = -6612.9375
$ws.Range("N134").Value = -9430.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16: H16=1310.1666, I16=1172.2, J16=2000, K16=1172.2, L16=2000, M16=-885.2, N16=-2574
$ws.Range("H16").Value = 1310.1666
$ws.Range("I16").Value = 1172.2
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1172.2
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -885.2
$ws.Range("N16").Value = -2574
# Row 58: H58=15648.8, I58=1356, J58=32621.5, K58=1356, L58=32621.5, M58=-1153, N58=-33027.5
$ws.Range("H58").Value = 15648.8
$ws.Range("I58").Value = 1356
$ws.Range("J58").Value = 32621.5
$ws.Range("K58").Value = 1356
$ws.Range("L58").Value = 32621.5
$ws.Range("M58").Value = -1153
$ws.Range("N58").Value = -33027.5
# Row 86: H86=12806.1, I86=2003.5, J86=15506.75, K86=2003.5, L86=15506.75, M86=-880.5, N86=-17752.75
$ws.Range("H86").Value = 12806.1
$ws.Range("I86").Value = 2003.5
$ws.Range("J86").Value = 15506.75
$ws.Range("K86").Value = 2003.5
$ws.Range("L86").Value = 15506.75
$ws.Range("M86").Value = -880.5
$ws.Range("N86").Value = -17752.75
# Row 89: H89=12806.1, I89=2003.5, J89=15506.75, K89=10017.5, L89=77533.75, M89=-4401.5, N89=-88765.75
$ws.Range("H89").Value = 12806.1
$ws.Range("I89").Value = 2003.5
$ws.Range("J89").Value = 15506.75
$ws.Range("K89").Value = 10017.5
$ws.Range("L89").Value = 77533.75
$ws.Range("M89").Value = -4401.5
$ws.Range("N89").Value = -88765.75
# Row 105: H105=1450, I105=933.3333, K105=933.3333, M105=813.6667
$ws.Range("H105").Value = 1450
$ws.Range("I105").Value = 933.3333
$ws.Range("K105").Value = 933.3333
$ws.Range("M105").Value = 813.6667
# Row 107: H107=1072.25, I107=431.57144, J107=1417.2307, K107=431.57144, L107=1417.2307, M107=1488.42856, N107=-5257.2307
$ws.Range("H107").Value = 1072.25
$ws.Range("I107").Value = 431.57144
$ws.Range("J107").Value = 1417.2307
$ws.Range("K107").Value = 431.57144
$ws.Range("L107").Value = 1417.2307
$ws.Range("M107").Value = 1488.42856
$ws.Range("N107").Value = -5257.2307
# Row 113: H113=1310.1666, I113=1172.2, J113=2000, K113=1172.2, L113=2000, M113=997.8, N113=-6340
$ws.Range("H113").Value = 1310.1666
$ws.Range("I113").Value = 1172.2
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1172.2
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 997.8
$ws.Range("N113").Value = -6340
# Row 122: H122=930.37036, I122=766.087, J122=1875, K122=2298.261, L122=5625, M122=151.739, N122=-10525
$ws.Range("H122").Value = 930.37036
$ws.Range("I122").Value = 766.087
$ws.Range("J122").Value = 1875
$ws.Range("K122").Value = 2298.261
$ws.Range("L122").Value = 5625
$ws.Range("M122").Value = 151.739
$ws.Range("N122").Value = -10525
# Row 124: H124=10796.893, I124=8738.565000000001, J124=20265.2, K124=8738.565000000001, L124=20265.2, M124=-6283.565000000001, N124=-25175.2
$ws.Range("H124").Value = 10796.893
$ws.Range("I124").Value = 8738.565000000001
$ws.Range("J124").Value = 20265.2
$ws.Range("K124").Value = 8738.565000000001
$ws.Range("L124").Value = 20265.2
$ws.Range("M124").Value = -6283.565000000001
$ws.Range("N124").Value = -25175.2
# Row 132: H132=2601.3667, I132=1955.6, J132=3892.9, K132=5866.799999999999, L132=11678.7, M132=-3336.799999999999, N132=-16738.7
$ws.Range("H132").Value = 2601.3667
$ws.Range("I132").Value = 1955.6
$ws.Range("J132").Value = 3892.9
$ws.Range("K132").Value = 5866.799999999999
$ws.Range("L132").Value = 11678.7
$ws.Range("M132").Value = -3336.799999999999
$ws.Range("N132").Value = -16738.7
# Row 136: H136=15648.8, I136=1356, J136=32621.5, K136=4068, L136=97864.5, M136=-1518, N136=-102964.5
$ws.Range("H136").Value = 15648.8
$ws.Range("I136").Value = 1356
$ws.Range("J136").Value = 32621.5
$ws.Range("K136").Value = 4068
$ws.Range("L136").Value = 97864.5
$ws.Range("M136").Value = -1518
$ws.Range("N136").Value = -102964.5
# Row 141: H141=74219.24000000001, J141=77103.375, L141=77103.375, N141=-87463.375
$ws.Range("H141").Value = 74219.24000000001
$ws.Range("J141").Value = 77103.375
$ws.Range("L141").Value = 77103.375
$ws.Range("N141").Value = -87463.375

$ws = $wb.Worksheets.Item("CUL")
# Row 115: H115=4062.7778, I115=1142.5, J115=6399, K115=3427.5, L115=19197, M115=-2252.5, N115=-21547
$ws.Range("H115").Value = 4062.7778
$ws.Range("I115").Value = 1142.5
$ws.Range("J115").Value = 6399
$ws.Range("K115").Value = 3427.5
$ws.Range("L115").Value = 19197
$ws.Range("M115").Value = -2252.5
$ws.Range("N115").Value = -21547
# Row 131: H131=730.7041, J131=760.8222, L131=2282.4666, N131=-12362.4666
$ws.Range("H131").Value = 730.7041
$ws.Range("J131").Value = 760.8222
$ws.Range("L131").Value = 2282.4666
$ws.Range("N131").Value = -12362.4666

$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80=20836790, I80=31253108, J80=4150, K80=31253108, L80=4150, M80=-31252110, N80=-6146
$ws.Range("H80").Value = 20836790
$ws.Range("I80").Value = 31253108
$ws.Range("J80").Value = 4150
$ws.Range("K80").Value = 31253108
$ws.Range("L80").Value = 4150
$ws.Range("M80").Value = -31252110
$ws.Range("N80").Value = -6146
# Row 83: H83=20836790, I83=31253108, J83=4150, K83=156265540, L83=20750, M83=-156260548, N83=-30734
$ws.Range("H83").Value = 20836790
$ws.Range("I83").Value = 31253108
$ws.Range("J83").Value = 4150
$ws.Range("K83").Value = 156265540
$ws.Range("L83").Value = 20750
$ws.Range("M83").Value = -156260548
$ws.Range("N83").Value = -30734
# Row 132: H132=20986.7, I132=4528.048, J132=59390.223, K132=13584.144, L132=178170.669, M132=-11054.144, N132=-183230.669
$ws.Range("H132").Value = 20986.7
$ws.Range("I132").Value = 4528.048
$ws.Range("J132").Value = 59390.223
$ws.Range("K132").Value = 13584.144
$ws.Range("L132").Value = 178170.669
$ws.Range("M132").Value = -11054.144
$ws.Range("N132").Value = -183230.669

$ws = $wb.Worksheets.Item("LTW")
# Row 61: H61=3594.45, I61=1954.9375, J61=10152.5, K61=1954.9375, L61=10152.5, M61=-1752.9375, N61=-10556.5
$ws.Range("H61").Value = 3594.45
$ws.Range("I61").Value = 1954.9375
$ws.Range("J61").Value = 10152.5
$ws.Range("K61").Value = 1954.9375
$ws.Range("L61").Value = 10152.5
$ws.Range("M61").Value = -1752.9375
$ws.Range("N61").Value = -10556.5
# Row 88: H88=33250, J88=33250, L88=33250, N88=-34106
$ws.Range("H88").Value = 33250
$ws.Range("J88").Value = 33250
$ws.Range("L88").Value = 33250
$ws.Range("N88").Value = -34106
# Row 91: H91=33250, J91=33250, L91=33250, N91=-36214
$ws.Range("H91").Value = 33250
$ws.Range("J91").Value = 33250
$ws.Range("L91").Value = 33250
$ws.Range("N91").Value = -36214
# Row 113: H113=3594.45, I113=1954.9375, J113=10152.5, K113=1954.9375, L113=10152.5, M113=215.0625, N113=-14492.5
$ws.Range("H113").Value = 3594.45
$ws.Range("I113").Value = 1954.9375
$ws.Range("J113").Value = 10152.5
$ws.Range("K113").Value = 1954.9375
$ws.Range("L113").Value = 10152.5
$ws.Range("M113").Value = 215.0625
$ws.Range("N113").Value = -14492.5
# Row 119: H119=33333.332, J119=33333.332, L119=33333.332, N119=-43009.332
$ws.Range("H119").Value = 33333.332
$ws.Range("J119").Value = 33333.332
$ws.Range("L119").Value = 33333.332
$ws.Range("N119").Value = -43009.332
# Row 122: H122=563064.4399999999, I122=728883.5600000001, K122=2186650.68, M122=-2184200.68
$ws.Range("H122").Value = 563064.4399999999
$ws.Range("I122").Value = 728883.5600000001
$ws.Range("K122").Value = 2186650.68
$ws.Range("M122").Value = -2184200.68
# Row 132: H132=243311.95, I132=327462.12, J132=3807.6155, K132=982386.36, L132=11422.8465, M132=-979856.36, N132=-16482.8465
$ws.Range("H132").Value = 243311.95
$ws.Range("I132").Value = 327462.12
$ws.Range("J132").Value = 3807.6155
$ws.Range("K132").Value = 982386.36
$ws.Range("L132").Value = 11422.8465
$ws.Range("M132").Value = -979856.36
$ws.Range("N132").Value = -16482.8465
# Row 136: H136=2103.5264, I136=1884.1428, K136=5652.428400000001, M136=-3102.428400000001
$ws.Range("H136").Value = 2103.5264
$ws.Range("I136").Value = 1884.1428
$ws.Range("K136").Value = 5652.428400000001
$ws.Range("M136").Value = -3102.428400000001

$ws = $wb.Worksheets.Item("WVR")
# Row 94: H94=20165, J94=20165, L94=20165, N94=-21967
$ws.Range("H94").Value = 20165
$ws.Range("J94").Value = 20165
$ws.Range("L94").Value = 20165
$ws.Range("N94").Value = -21967
# Row 100: H100=700, I100=750, K100=1500, M100=-959
$ws.Range("H100").Value = 700
$ws.Range("I100").Value = 750
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959
# Row 122: H122=1349, I122=1333.3334, J122=1372.5, K122=4000.0002, L122=4117.5, M122=-1550.0002, N122=-9017.5
$ws.Range("H122").Value = 1349
$ws.Range("I122").Value = 1333.3334
$ws.Range("J122").Value = 1372.5
$ws.Range("K122").Value = 4000.0002
$ws.Range("L122").Value = 4117.5
$ws.Range("M122").Value = -1550.0002
$ws.Range("N122").Value = -9017.5
# Row 132: H132=1099.2, I132=761.0476, K132=2283.1428, M132=246.8571999999999
$ws.Range("H132").Value = 1099.2
$ws.Range("I132").Value = 761.0476
$ws.Range("K132").Value = 2283.1428
$ws.Range("M132").Value = 246.8571999999999
# Row 136: H136=22941504, I136=30361776, J136=6118.636, K136=91085328, L136=18355.908, M136=-91082778, N136=-23455.908
$ws.Range("H136").Value = 22941504
$ws.Range("I136").Value = 30361776
$ws.Range("J136").Value = 6118.636
$ws.Range("K136").Value = 91085328
$ws.Range("L136").Value = 18355.908
$ws.Range("M136").Value = -91082778
$ws.Range("N136").Value = -23455.908

